$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.227.19'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '2.603.03'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '582.45'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = '142.83'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("E10").Value = '  -2.71%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '3.062.34'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.60'
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("D15").Value = '60.231.82'
$ws.Range("E15").Value = '  -1.29%  '
$ws.Range("E16").Value = '  -1.20%  '
$ws.Range("D17").Value = '2.607.40'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").Value = '11.33'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = '4.61'
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").Value = '347.29'
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Value = '6.89'
$ws.Range("E21").Value = '  -2.62%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").Value = '0.533'
$ws.Range("E23").Value = '  +3.50%  '
$ws.Range("D24").Value = '63.75'
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  +2.84%  '
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").Value = '0.0₃0797'
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("D30").Value = '169.12'
$ws.Range("E30").Value = '  +4.24%  '
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D34").Value = '1.33'
$ws.Range("E34").Value = '  +10.42%  '
$ws.Range("D35").Value = '4.24'
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("D36").Value = '0.991'
$ws.Range("E36").Value = '  +3.41%  '
$ws.Range("E37").Value = '  +2.47%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '316.44'
$ws.Range("E38").Value = '  +6.17%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '38.21'
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  +1.62%  '
$ws.Range("D41").Value = '0.845'
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '135.40'
$ws.Range("E42").Value = '  -2.93%  '
$ws.Range("D43").Value = '0.0993'
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").Value = '19.93'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").Value = '0.0548'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("D48").Value = '4.96'
$ws.Range("E48").Value = '  +5.83%  '
$ws.Range("D49").Value = '0.0241'
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("D51").Value = '10.73'
$ws.Range("E51").Value = '  +0.07%  '
